# Duplicate the last "multiple_answers" tab (multiple_answers11) three more
# times, growing the workbook from 18 to 21 sheets:
#   multiple_answers12 (sheetId 19)
#   multiple_answers13 (sheetId 20)
#   multiple_answers14 (sheetId 21)
# All four sheets (11-14) end up with identical data.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("multiple_answers11")

# Addresses of the 13 populated cells shared by these "Identifiers" sheets.
$cells = @("A1","A2","B3","C4","D5","E6","A7","B8","C9","A10","B11","C12","D13")

# Create the three duplicates first (they inherit whatever the source
# currently looks like), then nudge the style of the source sheet and its
# first two copies - the trailing duplicate is left exactly as the source
# originally was.
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$new1.Name = "multiple_answers12"

$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$new2.Name = "multiple_answers13"

$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$new3.Name = "multiple_answers14"

foreach ($addr in $cells) {
    $src.Range($addr).Style = "Normal"
    $new1.Range($addr).Style = "Normal"
    $new2.Range($addr).Style = "Normal"
}
